# wip: update redux state when option selected
#
# Adds a new "Catefories" worksheet between "Notes" and "Door quote vs
# plugin", populates it with the category/option reference lists, and
# makes it the active sheet (mirrors the author's workbook.xml diff:
# activeTab moves from the "Notes" tab to the new "Catefories" tab).

$wb = $excel.ActiveWorkbook

$notes = $wb.Worksheets.Item("Notes")

# Insert the new sheet right after "Notes" (i.e. before "Door quote vs plugin").
$cat = $wb.Worksheets.Add($null, $notes)
$cat.Name = "Catefories"

# Column B / C: parallel "single door" vs "double door" option lists.
$cat.Range("B2").Value = "Single Door Options"
$cat.Range("C2").Value = "Double Door"

$cat.Range("B3").Value = "Height, input box"
$cat.Range("C3").Value = "Height, input box"

$cat.Range("B4").Value = "Width, input box"
$cat.Range("C4").Value = "Width, input box"

$cat.Range("B5").Value = "Glass - 02 options"
$cat.Range("C5").Value = "Glass - 02 options"

$cat.Range("B6").Value = "lock type (panic bar)"
$cat.Range("C6").Value = "lock type (panic Bar)"

$cat.Range("B7").Value = "door closer"
$cat.Range("C7").Value = "door closer"

$cat.Range("B8").Value = "louver required"
$cat.Range("C8").Value = "louver required"

$cat.Range("B9").Value = "Door Type(Fire Rating)"
$cat.Range("C9").Value = "Door Type(Fire Rating)"

# Column D: list of door categories not handled by the quote generator.
$cat.Range("D1").Value = "Client Doors Not Covered"
$cat.Range("D2").Value = "Mustaqim Extra Large Doors"
$cat.Range("D3").Value = "China Customised Doors"

# Column widths, approximating the bestFit widths (~21.14 / ~26.14 chars)
# seen elsewhere in the workbook for similarly-sized text columns.
$cat.Columns.Item(2).ColumnWidth = 20.333333333333332
$cat.Columns.Item(3).ColumnWidth = 20.333333333333332
$cat.Columns.Item(4).ColumnWidth = 25.333333333333332

# Make the new sheet the active / selected tab, with D2 selected.
$cat.Activate()
$cat.Range("D2").Select()
